# error catching for excel upload
# Adds a new data row (row 3) to Sheet1 mirroring the existing layout,
# including a mailto hyperlink on the email column, matching the row
# already present for the first record (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data.
$ws.Range("A3").Value = "fgdg"
$ws.Range("B3").Value = "fgd"
$ws.Range("C3").Value = "dfgg"
$ws.Range("D3").Value = "fgdg"
$ws.Range("E3").Value = "gg@hh.com"

# Turn the email into a mailto: hyperlink, same as the existing E2 cell.
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:gg@hh.com", "", "", "gg@hh.com")

# Hyperlinks.Add auto-applies the hyperlink look (blue/underline) to the
# cell; the new row should keep the plain/default formatting used by its
# sibling cells, so restore E3's style from a plain neighbour on the row.
$ws.Range("A3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Move the active selection onto the newly added cell.
$ws.Range("E3").Select() | Out-Null
